$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header row (D1:E1) - copy style from C1 then set values
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

# Build data array for B2:E67
$arr = New-Object 'object[,]' 66,4
$arr[0,0] = -0.3239797079110615
$arr[0,1] = -0.3093622831163451
$arr[0,2] = -0.2951056968044581
$arr[0,3] = -0.281659107639499
$arr[1,0] = 0.1981554373629975
$arr[1,1] = 0.2070764539413797
$arr[1,2] = 0.2169806086954304
$arr[1,3] = 0.2268966025482406
$arr[2,0] = 0.1314601668315417
$arr[2,1] = 0.1431646595838115
$arr[2,2] = 0.155722516842765
$arr[2,3] = 0.1682514883545979
$arr[3,0] = -0.08508118290404185
$arr[3,1] = -0.0758442733719054
$arr[3,2] = -0.0668888729714044
$arr[3,3] = -0.05875336195302717
$arr[4,0] = 0.1290175634373698
$arr[4,1] = 0.1446280657153557
$arr[4,2] = 0.1595896857660735
$arr[4,3] = 0.1732006853884808
$arr[5,0] = -0.4604117615427986
$arr[5,1] = -0.4512205850314532
$arr[5,2] = -0.4414007614573875
$arr[5,3] = -0.4315106509063877
$arr[6,0] = -0.2705179443644665
$arr[6,1] = -0.2584775454032165
$arr[6,2] = -0.2457522000419311
$arr[6,3] = -0.2329860986406218
$arr[7,0] = -0.4205313625346355
$arr[7,1] = -0.4105889929375005
$arr[7,2] = -0.3983383449224077
$arr[7,3] = -0.384842012908065
$arr[8,0] = 0.3442059006747205
$arr[8,1] = 0.3588729555905856
$arr[8,2] = 0.3726771545349703
$arr[8,3] = 0.3849568979756428
$arr[9,0] = -0.2473425804904836
$arr[9,1] = -0.2385017336592683
$arr[9,2] = -0.2295179370961047
$arr[9,3] = -0.2209315470608058
$arr[10,0] = -0.1080640545265616
$arr[10,1] = -0.08619714041918237
$arr[10,2] = -0.06637391040061855
$arr[10,3] = -0.04914059488430657
$arr[11,0] = -0.01183244678392301
$arr[11,1] = -0.0118996560469052
$arr[11,2] = -0.01148368831016465
$arr[11,3] = -0.01106332846367018
$arr[12,0] = 0.01218659616566261
$arr[12,1] = 0.05218774876008093
$arr[12,2] = 0.08836179577735748
$arr[12,3] = 0.1200420080997457
$arr[13,0] = -0.05830895963806679
$arr[13,1] = -0.01411693445400752
$arr[13,2] = 0.02468361566765102
$arr[13,3] = 0.05764278502131422
$arr[14,0] = 0.2854292822568879
$arr[14,1] = 0.3464168503908811
$arr[14,2] = 0.4009131703446194
$arr[14,3] = 0.4480763620985666
$arr[15,0] = 0.563643104320362
$arr[15,1] = 0.5981719936045718
$arr[15,2] = 0.6273586310495535
$arr[15,3] = 0.6507584892659005
$arr[16,0] = 0.04863495344191104
$arr[16,1] = 0.03604968763938152
$arr[16,2] = 0.02511451067269434
$arr[16,3] = 0.01523885231427832
$arr[17,0] = 0.382040754536043
$arr[17,1] = 0.3972573364645569
$arr[17,2] = 0.4106325551481754
$arr[17,3] = 0.4215373026184711
$arr[18,0] = 0.2360995887169145
$arr[18,1] = 0.2935978017111165
$arr[18,2] = 0.3440917469733034
$arr[18,3] = 0.3869824702177221
$arr[19,0] = 0.4625655148839916
$arr[19,1] = 0.519205053466896
$arr[19,2] = 0.567757011915586
$arr[19,3] = 0.607896248786639
$arr[20,0] = 0.3416725059848368
$arr[20,1] = 0.3776076824557239
$arr[20,2] = 0.408419425252609
$arr[20,3] = 0.4336515563863449
$arr[21,0] = -0.09163552809647665
$arr[21,1] = -0.06027967599379161
$arr[21,2] = -0.03384627223570166
$arr[21,3] = -0.01249861779688013
$arr[22,0] = 4.530226679196272
$arr[22,1] = 4.564317052094395
$arr[22,2] = 4.517582289717446
$arr[22,3] = 4.405037572896482
$arr[23,0] = 0.503412429321434
$arr[23,1] = 0.4732361303223824
$arr[23,2] = 0.4442732946126758
$arr[23,3] = 0.4181124664046063
$arr[24,0] = 0.3995170796106945
$arr[24,1] = 0.3793329866124788
$arr[24,2] = 0.3567937659096443
$arr[24,3] = 0.3337624645120744
$arr[25,0] = 0.3194240927227214
$arr[25,1] = 0.2933557710250015
$arr[25,2] = 0.2662879357136232
$arr[25,3] = 0.2395090071076503
$arr[26,0] = 1.099008734566795
$arr[26,1] = 1.069421889075762
$arr[26,2] = 1.037770234966153
$arr[26,3] = 1.006044280207047
$arr[27,0] = 5.842635170195364
$arr[27,1] = 5.420665640241038
$arr[27,2] = 5.009052313905957
$arr[27,3] = 4.61151564955435
$arr[28,0] = 0.9901153388927074
$arr[28,1] = 0.9469909350344954
$arr[28,2] = 0.9032583720412695
$arr[28,3] = 0.8621207883974346
$arr[29,0] = -0.1661487883585251
$arr[29,1] = -0.2215551915663329
$arr[29,2] = -0.2717875579200011
$arr[29,3] = -0.3162886660036299
$arr[30,0] = 0.8213669964016488
$arr[30,1] = 0.7849417705186705
$arr[30,2] = 0.7496821511448362
$arr[30,3] = 0.7164581395144655
$arr[31,0] = 0.931716483315402
$arr[31,1] = 0.9049151329377539
$arr[31,2] = 0.8792998988590441
$arr[31,3] = 0.8557018617419659
$arr[32,0] = -0.6374160022790598
$arr[32,1] = -0.665537397882463
$arr[32,2] = -0.6934284907431428
$arr[32,3] = -0.71926888293991
$arr[33,0] = 0.8229529797785092
$arr[33,1] = 0.8139125278147399
$arr[33,2] = 0.8060482353208822
$arr[33,3] = 0.7987835857983867
$arr[34,0] = 0.7854012717963574
$arr[34,1] = 0.7679010701018232
$arr[34,2] = 0.7530758646338375
$arr[34,3] = 0.740278527206222
$arr[35,0] = 0.7599720416278469
$arr[35,1] = 0.7405970643419434
$arr[35,2] = 0.7240655154073049
$arr[35,3] = 0.7097441638984989
$arr[36,0] = 0.736573011468253
$arr[36,1] = 0.7172636952775983
$arr[36,2] = 0.7004664914487126
$arr[36,3] = 0.6855172595776524
$arr[37,0] = 0.5855179170203317
$arr[37,1] = 0.581361276923613
$arr[37,2] = 0.5782276845449104
$arr[37,3] = 0.5755308312763834
$arr[38,0] = 0.7561448125044548
$arr[38,1] = 0.7530707326558662
$arr[38,2] = 0.7503382177282049
$arr[38,3] = 0.7474351305597444
$arr[39,0] = 0.5624201316739293
$arr[39,1] = 0.5544038253167535
$arr[39,2] = 0.5486012461902333
$arr[39,3] = 0.5442323940907038
$arr[40,0] = 0.7178678997850312
$arr[40,1] = 0.6903968288805641
$arr[40,2] = 0.6662705758033712
$arr[40,3] = 0.6449146292890382
$arr[41,0] = 0.7235929828138905
$arr[41,1] = 0.7086778126739838
$arr[41,2] = 0.6961188762145806
$arr[41,3] = 0.6852093248637483
$arr[42,0] = 0.682510958967441
$arr[42,1] = 0.6742917402476475
$arr[42,2] = 0.6679788543289849
$arr[42,3] = 0.6628624381678773
$arr[43,0] = 0.676895140647332
$arr[43,1] = 0.6582783396307987
$arr[43,2] = 0.6428730729455288
$arr[43,3] = 0.6300272185980111
$arr[44,0] = -1.258761971176776
$arr[44,1] = -1.262984666059046
$arr[44,2] = -1.266611114531243
$arr[44,3] = -1.269546803655227
$arr[45,0] = -0.973519410757964
$arr[45,1] = -0.9787528595606133
$arr[45,2] = -0.9831894878414708
$arr[45,3] = -0.9868684582793954
$arr[46,0] = -0.8661238370075555
$arr[46,1] = -0.8725351967226652
$arr[46,2] = -0.8772899630234877
$arr[46,3] = -0.880609669866393
$arr[47,0] = -0.6378250778421266
$arr[47,1] = -0.6414807974495064
$arr[47,2] = -0.6439300494567668
$arr[47,3] = -0.645439551992443
$arr[48,0] = -0.0466675888936337
$arr[48,1] = -0.04809784013177344
$arr[48,2] = -0.04909847020926886
$arr[48,3] = -0.05004615127271864
$arr[49,0] = -0.8589558841379439
$arr[49,1] = -0.8643473716963861
$arr[49,2] = -0.8683839232675914
$arr[49,3] = -0.8712142478188347
$arr[50,0] = -0.8589558841379439
$arr[50,1] = -0.8643473716963861
$arr[50,2] = -0.8683839232675914
$arr[50,3] = -0.8712142478188347
$arr[51,0] = -1.083177447783509
$arr[51,1] = -1.096695474552693
$arr[51,2] = -1.107762864868827
$arr[51,3] = -1.116552343882033
$arr[52,0] = -0.1866984525573208
$arr[52,1] = -0.1862255579559371
$arr[52,2] = -0.1851842601744765
$arr[52,3] = -0.1839541820977837
$arr[53,0] = -0.993625540527033
$arr[53,1] = -0.9981244871884227
$arr[53,2] = -1.001981158626027
$arr[53,3] = -1.005214899595537
$arr[54,0] = -0.8870497724212896
$arr[54,1] = -0.8834438227905759
$arr[54,2] = -0.8808917751548818
$arr[54,3] = -0.8792413775838478
$arr[55,0] = -0.9455164010686488
$arr[55,1] = -0.9329141727921776
$arr[55,2] = -0.9217743867017627
$arr[55,3] = -0.9121361636797354
$arr[56,0] = -1.149612609280108
$arr[56,1] = -1.126906824887453
$arr[56,2] = -1.107034237706804
$arr[56,3] = -1.089863134883484
$arr[57,0] = -0.8587207998944131
$arr[57,1] = -0.8468034672072153
$arr[57,2] = -0.8356162180197534
$arr[57,3] = -0.8252526576901817
$arr[58,0] = -0.5160229371395958
$arr[58,1] = -0.4960781952889813
$arr[58,2] = -0.4779129566207586
$arr[58,3] = -0.4617633098022419
$arr[59,0] = 0.3697546559993954
$arr[59,1] = 0.3717217949089451
$arr[59,2] = 0.3742133307001154
$arr[59,3] = 0.3765940295057355
$arr[60,0] = -1.225893935686255
$arr[60,1] = -1.206867720350009
$arr[60,2] = -1.190202669450465
$arr[60,3] = -1.175841793752575
$arr[61,0] = -0.7629213430980695
$arr[61,1] = -0.7371004760351398
$arr[61,2] = -0.7125073593887912
$arr[61,3] = -0.6894967394457916
$arr[62,0] = -0.9064074911482318
$arr[62,1] = -0.9007103462505319
$arr[62,2] = -0.8937566734102732
$arr[62,3] = -0.8860632441221484
$arr[63,0] = -0.1290046537481677
$arr[63,1] = -0.109955290616581
$arr[63,2] = -0.09262663605800676
$arr[63,3] = -0.07738443671948579
$arr[64,0] = -0.8121577785326778
$arr[64,1] = -0.7887724666920192
$arr[64,2] = -0.7689191126246019
$arr[64,3] = -0.7524966189498451
$arr[65,0] = -0.7927760800172408
$arr[65,1] = -0.7605940848606052
$arr[65,2] = -0.7342722390111556
$arr[65,3] = -0.7133354046407581

$ws.Range("B2:E67").Value = $arr